# Update the cryptos price list with the latest scraped values.
# Columns: D = Price (text), E = Volume(1h) change (text, padded with spaces)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.505.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.843.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.020"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4579"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3871"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.18"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07896"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9931"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.81%  "
$ws.Range("E12").Value = "  -1.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.872.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.919"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.156"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06733"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001034"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.015"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.500.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.399"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.310"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.108"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.390"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9663"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09363"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.628"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.277"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.318"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02218"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05971"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.296"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.186"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.014"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5883"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1857"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.90%  "
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5553"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.890"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06689"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.045"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.016"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.08%  "
